# Apply the "Updated sizing results excel doc" edit.
#
# Summary of the change (from the supplied diff):
#  - The "Cold Case" sheet (2nd sheet) gets three new data rows (rows 3-5)
#    describing additional thermal-coating materials, with associated
#    Temp Min / Temp Max values in E1/F1 and updated column headers in
#    C2/D2 ("Temp Min [C]" / "Temp Max [C]" instead of "Binder"/"Color").
#  - The "Cold Case" sheet becomes the active/selected sheet (instead of
#    "Hot Case"), with a new selected cell.
#  - One of the new cells (the "-" placeholder company values) is entered
#    with a leading apostrophe (quote-prefixed text) the way a user would
#    type a value that Excel might otherwise try to interpret.
#  - Columns A-D on the Cold Case sheet get auto-fit to their new content.

$wb = $excel.ActiveWorkbook
$hotCase = $wb.Worksheets.Item("Hot Case")
$coldCase = $wb.Worksheets.Item("Cold Case")

# --- Update the header row on the Cold Case sheet -------------------------
# Eps eta values (row 1, columns E/F) next to the "Sierra Passive Louver"
# label.
$coldCase.Range("E1").Value = 0.13
$coldCase.Range("F1").Value = 0.14000000000000001

# --- New material rows ------------------------------------------------------
# (Entered in this order so newly-introduced shared strings land in the
# same order as the source workbook.)
$coldCase.Range("B3").Value = "NASA GSFC"
$coldCase.Range("C2").Value = "Temp Min [C]"
$coldCase.Range("D2").Value = "Temp Max [C]"
$coldCase.Range("A3").Value = "Dark Mirror SiO-Cr-Al"
$coldCase.Range("B4").Value = "'-"
$coldCase.Range("A5").Value = "Teflon Impregnated Anodized Titanium"
$coldCase.Range("B5").Value = "'-"
$coldCase.Range("A4").Value = "SiOx/VDA/0.5mil Kapton"

$coldCase.Range("E3").Value = 0.86
$coldCase.Range("F3").Value = 0.04

$coldCase.Range("E4").Value = 0.19
$coldCase.Range("F4").Value = 0.14000000000000001

$coldCase.Range("E5").Value = 0.76
$coldCase.Range("F5").Value = 0.26

# --- Column widths now that the sheet has real content in A:D -------------
$coldCase.Columns.Item(1).AutoFit() | Out-Null
$coldCase.Columns.Item(2).AutoFit() | Out-Null
$coldCase.Columns.Item(3).AutoFit() | Out-Null
$coldCase.Columns.Item(4).AutoFit() | Out-Null

# --- View/selection state ---------------------------------------------------
# Hot Case loses the "tabSelected" / scroll position it used to have...
$hotCase.Range("E9").Select() | Out-Null
# ...and Cold Case becomes the active sheet with F5 selected.
$coldCase.Activate()
$coldCase.Range("F5").Select() | Out-Null
